$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.553.34"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.739.48"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4923"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2677"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06293"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").Value = "1.735.39"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07050"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6149"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.589"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007376"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.87%  "
$ws.Range("D18").Value = "26.544.64"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9995"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").Value = "1.956.59"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.590"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.727"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.253"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "108.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.045"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08080"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.729"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04607"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.609"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.015"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6378"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8962"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.015"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.403"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.005"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01503"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.402"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3926"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.901"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05398"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.828"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.269"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.76"
$ws.Range("D51").Style = "Normal"
